$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (this also updates the Print_Area defined names
# that reference the sheet by name).
$ws.Name = "SAT Data"
$ws.PageSetup.PrintArea = "`$A`$1:`$U`$60"

# Change the style of Q6 to a new right-aligned, integer-format, boxed
# style (equivalent to the style already used by B6 but with a different
# border), then set the view so that the selection / scroll position
# match the saved state.
$ws.Range("Q6").NumberFormat = "0"
$ws.Range("Q6").Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$ws.Range("Q6").Borders.Item(7).Weight = 2      # xlThin
$ws.Range("Q6").Borders.Item(8).LineStyle = 1   # xlEdgeTop
$ws.Range("Q6").Borders.Item(8).Weight = 2      # xlThin
$ws.Range("Q6").Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$ws.Range("Q6").Borders.Item(9).Weight = 2      # xlThin
$ws.Range("Q6").HorizontalAlignment = -4152     # xlRight
$ws.Range("Q6").VerticalAlignment = -4108       # xlCenter

$ws.Range("Q6").Select()
$ws.Application.ActiveWindow.ScrollColumn = 8   # Column H
$ws.Application.ActiveWindow.ScrollRow = 1
